# Weekly update: insert a new fruit/vegetable price record as row 3,
# pushing the existing historical rows (old 3..23) down to 4..24.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(3).Insert()

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Macroferia Regional de Talca"
$ws.Range("C3").Value = "Maule"
$ws.Range("D3").Value = 44462
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 100112026
$ws.Range("G3").Value = "Haba"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("N3").Value = "`$/saco 25 kilos"
$ws.Range("O3").Value = "Región de O'Higgins"
$ws.Range("P3").Value = 400
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
